$d = $word.ActiveDocument

# 1. Créditos-aula: 2 -> 4
$d.Content.Find.Execute("Créditos-aula: 2", $true, $false, $false, $false, $false, $true, 1, $false, "Créditos-aula: 4", 2)

# 2. Carga horária: 30 h -> 60 h
$d.Content.Find.Execute("Carga horária: 30 h", $true, $false, $false, $false, $false, $true, 1, $false, "Carga horária: 60 h", 2)

# 3. Ativação: 01/01/2012 -> 01/01/2025
$d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2)

# 4. Objetivos paragraph
$d.Content.Find.Execute("Apresentar e analisar os conceitos básicos de monitoramento, suas aplicações práticas e as interfaces com os demais instrumentos de Política Ambiental.", $true, $false, $false, $false, $false, $true, 1, $false, "Apresentar e analisar as técnicas de monitoramento dos principais parâmetros ambientais exigidos pelas legislações em vigor.", 2)

# 5. Docente(s) Responsável(eis) paragraph: swap/replace the two teacher lines.
#    Original runs (same formatting, both plain):
#      run A: "5840938 - Marcelo Rodrigues de Holanda" + <w:br/>
#      run B: "8855158 - Morun Bernardino Neto"
#    Target:
#      run A: "8855158 - Morun Bernardino Neto" + <w:br/>
#      run B: "7455355 - Robson da Silva Rocha"
#    Find the paragraph dynamically (don't hardcode char offsets, they shift with earlier edits).
$docentesPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*5840938 - Marcelo Rodrigues de Holanda*8855158 - Morun Bernardino Neto*") {
        $docentesPara = $p
    }
}
$pStart = $docentesPara.Range.Start
$pEnd = $docentesPara.Range.End
$nameA = "5840938 - Marcelo Rodrigues de Holanda"
$nameB = "8855158 - Morun Bernardino Neto"
$newA = "8855158 - Morun Bernardino Neto"
$newB = "7455355 - Robson da Silva Rocha"
$rangeB = $d.Range($pStart + $nameA.Length + 1, $pEnd)
$rangeB.Text = $newB
$rangeA = $d.Range($pStart, $pStart + $nameA.Length)
$rangeA.Text = $newA

# 6. Programa resumido paragraph
$d.Content.Find.Execute("Monitoramento da qualidade ambiental.", $true, $false, $false, $false, $false, $true, 1, $false, "Conceito de monitoramento. Amostragem. Técnicas alternativas para cada parâmetro a ser monitorado. Relação custo e aplicabilidade.", 2)

# 7. Programa paragraph
$d.Content.Find.Execute("Conceitos de qualidade ambiental, poluição, padrões de qualidade e de emissão. Conceito de monitoramento. Amostragem. Sistemas de monitoramento. Índices de qualidade. Monitoramento como parte integrante de sistema de gestão ambiental.", $true, $false, $false, $false, $false, $true, 1, $false, "- Conceito de monitoramento.- Técnicas de amostragens e suas especificidades para cada parâmetro.- Técnicas e equipamentos para monitorar: carga orgânica, sólidos, íons, metais, atividade biológica e outros parâmetros de importância ambiental.- Química Verde no monitoramento ambiental", 2)

# 8. Método
$d.Content.Find.Execute("Aula expositiva e exercícios dirigidos.", $true, $false, $false, $false, $false, $true, 1, $false, "Avaliação baseada em provas, exercícios, trabalhos práticos e relatórios.", 2)

# 9. Critério
$d.Content.Find.Execute("Média ponderada de exercícios e provas.", $true, $false, $false, $false, $false, $true, 1, $false, "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.", 2)

# 10. Norma de recuperação
$d.Content.Find.Execute("Prova única com nota igual ou superior a 5,0.", $true, $false, $false, $false, $false, $true, 1, $false, "1 (uma) prova de recuperação (R), sendo considerado aprovado se R >= 5,0.", 2)
